$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): 想去人数 (F) bumped for a few rows
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 541
$ws1.Range("F9").Value = 392
$ws1.Range("F10").Value = 3429
$ws1.Range("F11").Value = 42

# Sheet "演出" (index 2): 最低票价 (G) updated for row 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = 98

# Sheet "全部类型" (index 4): mirrors the same underlying events
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("G3").Value = 98
$ws4.Range("F4").Value = 541
$ws4.Range("F10").Value = 392
$ws4.Range("F11").Value = 3429
$ws4.Range("F12").Value = 42
